# Update column G ("K") values on Sheet1 per the regenerated save_data.
# Column G holds the strike-count metric ("K"); the values below were
# recalculated (std/mean, calc and write s_vals) and replace the
# previous Strike#-derived numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2 = 0
    3 = 4
    4 = 2
    5 = 4
    7 = 4
    8 = 2
    9 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
